$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04
$ws.Range("D2").Value = 0.005
$ws.Range("C3").Value = 0.04
$ws.Range("D3").Value = 0.005

$ws.Range("D3").Select()
